# The commit removes one duplicate leaderboard entry:
#   row 3 = 2024-02-02 timestamp, name "Anamika", roll "B23428", branch "ME"
# which is a duplicate of the "Anamika" / "b23428" entry already present a
# few rows down (now the correct, surviving record). Deleting the row
# shifts every later row up by one (old row 4 -> new row 3, ..., old row 57
# -> new row 56), shrinking the used range from A1:F57 to A1:F56 and
# dropping the now-unreferenced "B23428" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Delete()

# The one hyperlink in the sheet (the beecrowd profile link) lived on D39
# and therefore needs to end up on D38 after the shift above. This runtime
# doesn't re-target an existing Hyperlink's anchor range when rows move, so
# drop the (now stale) hyperlink and recreate it over the correct cell.
$target = "https://www.beecrowd.com.br/judge/en/profile/948533"
$ws.Range("D39").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D38"), $target)

# Reflect the author's final cursor position/selection after the edit
# (no more frozen/scrolled top-left cell, selection parked at D60).
$ws.Range("D60").Select()
